$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right after the header (pushes the existing
# rows 4-12 down to 5-13, carrying their formatting/style with them).
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44659
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112043
$ws.Range("G4").Value = "Pepino dulce"
$ws.Range("H4").Value = "Cultivar IV Región"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 861
$ws.Range("Q4").Value = 18
$ws.Range("R4").Value = "Hortaliza"
